$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (before) values of the rotating fields for rows 2-5
$rows = 2..5
$cols = @("A","B","E","F","G","H","I","J","Q","R")

$data = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $data[$r] = $rowData
}

# Rotate: new row2 = old row5, new row3 = old row2, new row4 = old row3, new row5 = old row4
$mapping = @{ 2 = 5; 3 = 2; 4 = 3; 5 = 4 }

foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $data[$srcRow][$c]
    }
}
